$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.253.85"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "'1.908.76"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'0.720"
$ws.Range("E5").Value = "  +8.31%  "
$ws.Range("D6").Value = "'255.97"
$ws.Range("E6").Value = "  +3.93%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'40.66"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("E9").Value = "  +6.96%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "'0.0761"
$ws.Range("E11").Value = "  +5.58%  "
$ws.Range("D12").Value = "'0.0987"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "'2.187.08"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "'12.84"
$ws.Range("E14").Value = "  +6.24%  "
$ws.Range("D15").Value = "'0.727"
$ws.Range("E15").Value = "  +3.87%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'1.966.96"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'4.96"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").Value = "'35.257.00"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "'74.67"
$ws.Range("E19").Value = "  +3.31%  "
$ws.Range("D20").Value = "'0.0₃0849"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").Value = "'243.62"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").Value = "'13.04"
$ws.Range("E22").Value = "  +4.09%  "
$ws.Range("E23").Value = "  +5.48%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'2.44"
$ws.Range("E25").Value = "  +6.38%  "
$ws.Range("D26").Value = "'2.44"
$ws.Range("E26").Value = "  +4.25%  "
$ws.Range("D27").Value = "'166.38"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D28").Value = "'8.68"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("D29").Value = "'18.74"
$ws.Range("E29").Value = "  +1.58%  "
$ws.Range("E30").Value = "  +3.44%  "
$ws.Range("D31").Value = "'4.128.95"
$ws.Range("E31").Value = "  +19.46%  "
$ws.Range("E32").Value = "  +5.69%  "
$ws.Range("E33").Value = "  +14.27%  "
$ws.Range("E34").Value = "  +22.12%  "
$ws.Range("E35").Value = "  +3.62%  "
$ws.Range("D36").Value = "'4.23"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").Value = "'0.912"
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").Value = "'17.23"
$ws.Range("E40").Value = "  +5.20%  "
$ws.Range("D41").Value = "'0.0218"
$ws.Range("E41").Value = "  +4.51%  "
$ws.Range("D42").Value = "'96.52"
$ws.Range("E42").Value = "  +7.15%  "
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("D44").Value = "'0.0654"
$ws.Range("E44").Value = "  +2.94%  "
$ws.Range("D45").Value = "'1.334.62"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").Value = "'2.42"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D48").Value = "'6.72"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").Value = "'44.91"
$ws.Range("E50").Value = "  -6.62%  "
$ws.Range("E51").Value = "  +6.28%  "
